# Atualizei dados bibi e add
# - Corrige 3 valores de total_venda no bloco de agosto/2025 (linhas 5, 7 e 10)
# - Insere um novo registro (dia 14, agosto/2025) na linha 11, empurrando
#   todas as linhas seguintes uma posicao para baixo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige valores existentes
$ws.Range("B5").Value = 20255.27
$ws.Range("B7").Value = 13698.08
$ws.Range("B10").Value = 24144.7

# Insere uma nova linha na posicao 11 (desloca as demais para baixo)
$ws.Rows.Item(11).Insert()

# Preenche a nova linha com o novo registro diario
$ws.Range("A11").Value = 14
$ws.Range("B11").Value = 14956.89
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 2025
$ws.Range("E11").Value = "08/2025"
